$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-obsolete trailing rows (old table had rows through 22; new table ends at row 14)
$ws.Rows("15:22").Delete()

# Rewrite rows 2-14 with the updated year range (2010-2022) and figures
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 308600
$ws.Range("C2").Value = 11473000
$ws.Range("D2").Value = 10882000

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 168700
$ws.Range("C3").Value = 12399000
$ws.Range("D3").Value = 11769800

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 157800
$ws.Range("C4").Value = 12107000
$ws.Range("D4").Value = 11329400

$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 113400
$ws.Range("C5").Value = 12391000
$ws.Range("D5").Value = 11872100

$ws.Range("A6").Value = "2014年"
$ws.Range("B6").Value = 14400
$ws.Range("C6").Value = 12850000
$ws.Range("D6").Value = 11970500

$ws.Range("A7").Value = "2015年"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 13558000
$ws.Range("D7").Value = 12626700

$ws.Range("A8").Value = "2016年"
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 13371000
$ws.Range("D8").Value = 12600100

$ws.Range("A9").Value = "2017年"
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = 13632000
$ws.Range("D9").Value = 13103500

$ws.Range("A10").Value = "2018年"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 13831000
$ws.Range("D10").Value = 13496600

$ws.Range("A11").Value = "2019年"
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 14122000
$ws.Range("D11").Value = 13813500

$ws.Range("A12").Value = "2020年"
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 14934000
$ws.Range("D12").Value = 14436900

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 18096000
$ws.Range("D13").Value = 17348300

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 18913241
$ws.Range("D14").ClearContents()
